# Insert a new weekly record at row 83 (Macroferia Regional de Talca - Mango).
# This pushes the previous rows 83..115 down to 84..116, growing the used
# range from A1:T115 to A1:T116.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(83).Insert()

$ws.Cells.Item(83, 1).Value = 5
$ws.Cells.Item(83, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(83, 3).Value = "Maule"
$ws.Cells.Item(83, 4).Value = 44627
$ws.Cells.Item(83, 5).Value = 7
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100108
$ws.Cells.Item(83, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(83, 9).Value = 100108002
$ws.Cells.Item(83, 10).Value = "Mango"
$ws.Cells.Item(83, 11).Value = "Sin especificar"
$ws.Cells.Item(83, 12).Value = "Primera"
$ws.Cells.Item(83, 13).Value = 200
$ws.Cells.Item(83, 14).Value = 7000
$ws.Cells.Item(83, 15).Value = 7000
$ws.Cells.Item(83, 16).Value = 7000
$ws.Cells.Item(83, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(83, 18).Value = "Ecuador"
$ws.Cells.Item(83, 19).Value = 1750
$ws.Cells.Item(83, 20).Value = 4
